# Updated cryptos list on Tue Mar 19 04:48:22 UTC 2024 with GitHub Actions
# Applies the price/volume refresh (and the two name re-orderings) to the
# cryptocurrency listing on Sheet1. Column B = Coin, C = Link, D = Price,
# E = Volume(1h). Only the cells whose contents actually changed are
# touched; the leading apostrophe forces each D/E write to stay a text
# value (matching the original formatting) instead of being
# auto-coerced into a number/percentage by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="'65.536.14"; E="'  -4.26%  "},
    @{Row=3; D="'3.392.28"; E="'  -6.38%  "},
    @{Row=4; E="'  +0.10%  "},
    @{Row=5; D="'186.67"; E="'  -8.61%  "},
    @{Row=6; D="'531.12"; E="'  -6.61%  "},
    @{Row=7; D="'0.609"; E="'  -2.19%  "},
    @{Row=8; D="'3.385.63"; E="'  -6.43%  "},
    @{Row=9; E="'  -0.03%  "},
    @{Row=10; D="'0.631"; E="'  -6.91%  "},
    @{Row=11; D="'59.24"; E="'  -3.59%  "},
    @{Row=12; D="'0.134"; E="'  -12.07%  "},
    @{Row=13; D="'0.0000257"; E="'  -11.11%  "},
    @{Row=14; D="'9.32"; E="'  -7.78%  "},
    @{Row=15; D="'3.941.77"; E="'  -5.99%  "},
    @{Row=16; E="'  -2.83%  "},
    @{Row=17; D="'3.400.14"; E="'  -5.95%  "},
    @{Row=18; D="'65.229.90"; E="'  -4.40%  "},
    @{Row=19; D="'17.53"; E="'  -8.14%  "},
    @{Row=20; D="'11.21"; E="'  -9.66%  "},
    @{Row=21; D="'0.979"; E="'  -9.26%  "},
    @{Row=22; D="'375.16"; E="'  -7.43%  "},
    @{Row=23; D="'82.30"; E="'  -3.82%  "},
    @{Row=24; D="'3.76"; E="'  -10.11%  "},
    @{Row=25; D="'10.89"; E="'  -16.01%  "},
    @{Row=26; D="'3.72"; E="'  -4.99%  "},
    @{Row=27; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="'2.67"; E="'  -9.12%  "},
    @{Row=28; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="'11.51"; E="'  -8.86%  "},
    @{Row=29; D="'8.59"; E="'  -8.86%  "},
    @{Row=30; D="'692.33"; E="'  +3.01%  "},
    @{Row=31; D="'29.83"; E="'  -6.07%  "},
    @{Row=32; D="'6.77"; E="'  -17.28%  "},
    @{Row=33; D="'11.25"; E="'  -8.51%  "},
    @{Row=34; D="'61.26"; E="'  -4.30%  "},
    @{Row=35; E="'  -7.10%  "},
    @{Row=36; E="'  -0.12%  "},
    @{Row=37; D="'36.74"; E="'  -13.21%  "},
    @{Row=38; D="'0.385"; E="'  -9.15%  "},
    @{Row=39; E="'  -0.01%  "},
    @{Row=40; E="'  -5.77%  "},
    @{Row=41; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="'2.871.68"; E="'  -12.26%  "},
    @{Row=42; B="ThetaToken"; C="https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"; D="'2.79"; E="'  -12.34%  "},
    @{Row=43; B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="'2.68"; E="'  -3.96%  "},
    @{Row=44; D="'0.0401"; E="'  -4.49%  "},
    @{Row=45; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="'0.0₃0626"; E="'  -19.43%  "},
    @{Row=46; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="'2.37"; E="'  -14.30%  "},
    @{Row=47; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="'0.127"; E="'  -3.92%  "},
    @{Row=48; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="'137.78"; E="'  -1.13%  "},
    @{Row=49; B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="'2.65"; E="'  -3.56%  "},
    @{Row=50; D="'2.85"; E="'  -8.53%  "},
    @{Row=51; B="THORChain"; C="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D="'7.72"; E="'  -13.30%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
